# Apply "custom accuracy" rounding to row 5 (reduce to 2 decimal places)
# and remove row 6 (delete the last data row), per commit:
#   "custom accuracy + 데이터 1000개"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: overwrite row 5 values with their "custom accuracy" (2 decimal place) values ---
$row5 = @{
    "B5"  = 3.42
    "C5"  = 2.27
    "D5"  = 0.01
    "E5"  = 5.6
    "F5"  = 5
    "G5"  = 2.55
    "H5"  = 13.16
    "I5"  = 3.11
    "J5"  = 1.6
    "K5"  = 2.74
    "L5"  = 2.37
    "M5"  = 2.27
    "N5"  = 0.74
    "O5"  = 2.02
    "P5"  = 3.56
    "Q5"  = 1.72
    "R5"  = 0.26
    "S5"  = 0
    "T5"  = 26.87
    "U5"  = 6.53
    "V5"  = 2.25
    "W5"  = 4.54
    "X5"  = 2.12
    "Y5"  = 0.3
    "Z5"  = 5.99
    "AA5" = 1.82
    "AB5" = 1.44
    "AC5" = 1.74
    "AD5" = 2.96
    "AE5" = 0.52
    "AF5" = 11.91
    "AG5" = 1.03
    "AH5" = 2.44
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# --- Step 2: delete row 6 entirely (shifts rows up, shrinking the used range) ---
$ws.Rows.Item(6).Delete()

$wb.Save()
